$wb = $excel.ActiveWorkbook

# --- Controls sheet: n_sims 200 -> 100 ---
$wsControls = $wb.Worksheets.Item("Controls")
$wsControls.Range("B2").Value = 100

# --- Growth_Param sheet: update L_inf / k values, fix r0 bug values ---
$wsGrowth = $wb.Worksheets.Item("Growth_Param")
$wsGrowth.Range("B2").Value = 0.2
$wsGrowth.Range("A3").Value = 95
$wsGrowth.Range("B3").Value = 70

# --- Selex sheet: move selection to B4 ---
$wsSelex = $wb.Worksheets.Item("Selex")
$wsSelex.Range("B4").Select()

# --- Make Growth_Param the active/selected sheet with selection on B4 ---
# (done last so it ends up as the active tab/window selection)
$wsGrowth.Activate()
$wsGrowth.Range("B4").Select()
